# Auto-generated Excel COM-interop script
# Applies the "Kraken_Profits" market-price refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

$ALC_updates = @{
    "H17" = 1750
    "J17" = 1900
    "L17" = 5700
    "N17" = -6036
    "H19" = 599
    "I19" = 798
    "J19" = 499.5
    "K19" = 798
    "L19" = 499.5
    "M19" = -623
    "N19" = -849.5
    "H33" = 408.82352
    "I33" = 126.666664
    "K33" = 126.666664
    "M33" = 102.333336
    "H38" = 29.857143
    "I38" = 29.857143
    "J38" = 0
    "K38" = 89.57142899999999
    "L38" = 0
    "M38" = 282.428571
    "H39" = 12.75
    "I39" = 12.75
    "K39" = 38.25
    "M39" = 257.75
    "H53" = 163
    "I53" = 129
    "J53" = 214
    "K53" = 129
    "L53" = 214
    "M53" = 508
    "N53" = -1488
    "H70" = 0
    "I70" = 0
    "J70" = 0
    "K70" = 0
    "L70" = 0
    "H73" = 0
    "I73" = 0
    "J73" = 0
    "K73" = 0
    "L73" = 0
    "H137" = 2055.875
    "I137" = 1991.3334
    "J137" = 2249.5
    "K137" = 5974.0002
    "L137" = 6748.5
    "M137" = -3424.0002
    "N137" = -11848.5
}
foreach ($cell in $ALC_updates.Keys) {
    $ws.Range($cell).Value = $ALC_updates[$cell]
}

$ALC_clears = @("N38", "M70", "N70", "M73", "N73")
foreach ($cell in $ALC_clears) {
    $ws.Range($cell).ClearContents()
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

$ARM_updates = @{
    "H46" = 20050.666
    "I46" = 20000
    "J46" = 20076
    "K46" = 20000
    "L46" = 20076
    "M46" = -19681
    "N46" = -20714
    "H61" = 3577.5715
    "I61" = 3007.75
    "J61" = 6996.5
    "K61" = 3007.75
    "L61" = 6996.5
    "M61" = -2795.75
    "N61" = -7420.5
    "H88" = 1425
    "J88" = 873.5
    "L88" = 873.5
    "N88" = -1685.5
    "H91" = 1425
    "J91" = 873.5
    "L91" = 873.5
    "N91" = -3681.5
    "I97" = 1240.3334
    "J97" = 4350
    "K97" = 1240.3334
    "L97" = 4350
    "M97" = -744.3334
    "N97" = -5342
    "H136" = 3577.5715
    "I136" = 3007.75
    "J136" = 6996.5
    "K136" = 9023.25
    "L136" = 20989.5
    "M136" = -6473.25
    "N136" = -26089.5
}
foreach ($cell in $ARM_updates.Keys) {
    $ws.Range($cell).Value = $ARM_updates[$cell]
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

$BSM_updates = @{
    "H82" = 11164.5
    "I82" = 11164.5
    "K82" = 11164.5
    "M82" = -10781.5
    "H85" = 11164.5
    "I85" = 11164.5
    "K85" = 11164.5
    "M85" = -9838.5
    "H97" = 10262.75
    "I97" = 10262.75
    "K97" = 10262.75
    "M97" = -9271.75
}
foreach ($cell in $BSM_updates.Keys) {
    $ws.Range($cell).Value = $BSM_updates[$cell]
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

$CRP_updates = @{
    "H22" = 600
    "I22" = 600
    "K22" = 600
    "M22" = -250
    "H58" = 0
    "I58" = 0
    "K58" = 0
    "H62" = 4519.6
    "I62" = 4866.3335
    "J62" = 3999.5
    "K62" = 4866.3335
    "L62" = 3999.5
    "M62" = -4242.3335
    "N62" = -5247.5
    "H65" = 4519.6
    "I65" = 4866.3335
    "J65" = 3999.5
    "K65" = 24331.6675
    "L65" = 19997.5
    "M65" = -21211.6675
    "N65" = -26237.5
    "H136" = 0
    "I136" = 0
    "K136" = 0
}
foreach ($cell in $CRP_updates.Keys) {
    $ws.Range($cell).Value = $CRP_updates[$cell]
}

$CRP_clears = @("M58", "M136")
foreach ($cell in $CRP_clears) {
    $ws.Range($cell).ClearContents()
}

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

$CUL_updates = @{
    "H2" = 25.428572
    "I2" = 15.428572
    "J2" = 35.42857
    "K2" = 92.571432
    "L2" = 212.57142
    "M2" = 20.428568
    "N2" = -438.57142
    "H11" = 525.1667
    "I11" = 525.1667
    "K11" = 1575.5001
    "M11" = -1435.5001
    "H113" = 887.5
}
foreach ($cell in $CUL_updates.Keys) {
    $ws.Range($cell).Value = $CUL_updates[$cell]
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

$GSM_updates = @{
    "H107" = 441.66666
    "I107" = 329.4
    "J107" = 1003
    "K107" = 329.4
    "M107" = 1590.6
    "N107" = -4843
    "H132" = 4500.2856
    "I132" = 3791.5
    "K132" = 11374.5
    "M132" = -8844.5
}
foreach ($cell in $GSM_updates.Keys) {
    $ws.Range($cell).Value = $GSM_updates[$cell]
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

$LTW_updates = @{
    "H22" = 977.7778
    "I22" = 971.4286
    "J22" = 1000
    "K22" = 971.4286
    "L22" = 1000
    "M22" = -676.4286
    "N22" = -1590
    "H27" = 977.7778
    "I27" = 971.4286
    "J27" = 1000
    "K27" = 971.4286
    "L27" = 1000
    "M27" = -864.4286
    "N27" = -1214
}
foreach ($cell in $LTW_updates.Keys) {
    $ws.Range($cell).Value = $LTW_updates[$cell]
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

$WVR_updates = @{
    "H96" = 0
    "I96" = 0
    "K96" = 0
    "H107" = 4500
    "I107" = 2000
    "K107" = 6000
    "M107" = -4080
    "H132" = 2786.625
    "I132" = 1049
    "K132" = 3147
    "M132" = -617
}
foreach ($cell in $WVR_updates.Keys) {
    $ws.Range($cell).Value = $WVR_updates[$cell]
}

$WVR_clears = @("M96")
foreach ($cell in $WVR_clears) {
    $ws.Range($cell).ClearContents()
}
